# Apply the "10031723" auto-commit edit:
#  - Append 5 new data rows (15-19) to the Report sheet, cloning the
#    alternating row-striping format from the existing template rows
#    (row 13 = odd-row style, row 14 = even-row style).
#  - Mark the "工作內容" (AC) / "報修說明" (P) cells that should wrap text,
#    matching rows 13/14's pattern.
#  - Extend the print area to the new bottom row and update the active
#    selection, mirroring what Excel does automatically when a user keys
#    in new rows and the workbook is re-saved.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Row 14's P/AC cells pick up the "wrap" variant of their style ---
$ws.Range("P14").WrapText = $true
$ws.Range("AC14").WrapText = $true

# --- 2. Clone formatting for the five new rows from the matching template row ---
# Odd data rows (15, 17, 19) mirror row 13's style band (fill + borders).
# Even data rows (16, 18) mirror row 14's style band.
$ws.Range("A13:AK13").Copy() | Out-Null
$ws.Range("A15:AK15").PasteSpecial(-4122) | Out-Null

$ws.Range("A14:AK14").Copy() | Out-Null
$ws.Range("A16:AK16").PasteSpecial(-4122) | Out-Null

$ws.Range("A13:AK13").Copy() | Out-Null
$ws.Range("A17:AK17").PasteSpecial(-4122) | Out-Null

$ws.Range("A14:AK14").Copy() | Out-Null
$ws.Range("A18:AK18").PasteSpecial(-4122) | Out-Null

$ws.Range("A13:AK13").Copy() | Out-Null
$ws.Range("A19:AK19").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# The "報修說明" / "工作內容" cells in rows 15-18 wrap (row 19 keeps the
# plain, non-wrapped variant, matching the source workbook).
$ws.Range("P15").WrapText = $true
$ws.Range("AC15").WrapText = $true
$ws.Range("P16").WrapText = $true
$ws.Range("AC16").WrapText = $true
$ws.Range("P17").WrapText = $true
$ws.Range("AC17").WrapText = $true
$ws.Range("P18").WrapText = $true
$ws.Range("AC18").WrapText = $true

# --- 3. Row 15 (item 13) - 三重仁義店 / THILF04241 ---
$ws.Range("A15").Value = 13
$ws.Range("B15").Value = "服務"
$ws.Range("C15").Value = 2025100662
$ws.Range("F15").Value = 4241
$ws.Range("G15").Value = "三重仁義店"
$ws.Range("H15").Value = "新北市三重區"
$ws.Range("Q15").Value = "THILF04241"
$ws.Range("R15").Value = "新北一"
$ws.Range("S15").Value = "吳宗鴻"
$ws.Range("T15").Value = 1
$ws.Range("U15").Value = "已完工"
$ws.Range("V15").Value = "2025-10-03 13:58:34"
$ws.Range("W15").Value = "2025-10-03 13:30:00"
$ws.Range("X15").Value = "2025-10-03 13:50:00"
$ws.Range("Z15").Value = 0.3
$ws.Range("AB15").Value = "到場處理"
$ws.Range("AC15").Value = "PMQ4"
$ws.Range("AD15").Value = "O"
$ws.Range("AK15").Value = "O"

# --- 4. Row 16 (item 14) - 三重公園店 / THILF04352 ---
$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "服務"
$ws.Range("C16").Value = 2025100682
$ws.Range("F16").Value = 4352
$ws.Range("G16").Value = "三重公園店"
$ws.Range("H16").Value = "新北市三重區"
$ws.Range("Q16").Value = "THILF04352"
$ws.Range("R16").Value = "新北一"
$ws.Range("S16").Value = "吳宗鴻"
$ws.Range("T16").Value = 1
$ws.Range("U16").Value = "已完工"
$ws.Range("V16").Value = "2025-10-03 14:34:21"
$ws.Range("W16").Value = "2025-10-03 14:10:00"
$ws.Range("X16").Value = "2025-10-03 14:30:00"
$ws.Range("Z16").Value = 0.3
$ws.Range("AB16").Value = "到場處理"
$ws.Range("AC16").Value = "PMQ4"
$ws.Range("AD16").Value = "O"
$ws.Range("AK16").Value = "O"

# --- 5. Row 17 (item 15) - 北縣天龍店 / THILF03840 ---
$ws.Range("A17").Value = 15
$ws.Range("B17").Value = "服務"
$ws.Range("C17").Value = 2025100705
$ws.Range("F17").Value = 3840
$ws.Range("G17").Value = "北縣天龍店"
$ws.Range("H17").Value = "新北市三重區"
$ws.Range("Q17").Value = "THILF03840"
$ws.Range("R17").Value = "新北一"
$ws.Range("S17").Value = "吳宗鴻"
$ws.Range("T17").Value = 1
$ws.Range("U17").Value = "已完工"
$ws.Range("V17").Value = "2025-10-03 15:35:02"
$ws.Range("W17").Value = "2025-10-03 15:10:00"
$ws.Range("X17").Value = "2025-10-03 15:34:00"
$ws.Range("Z17").Value = 0.4
$ws.Range("AB17").Value = "到場處理"
$ws.Range("AC17").Value = "PMQ4"
$ws.Range("AD17").Value = "O"
$ws.Range("AK17").Value = "O"

# --- 6. Row 18 (item 16) - 三重仁旺店 / THILF04397 ---
$ws.Range("A18").Value = 16
$ws.Range("B18").Value = "服務"
$ws.Range("C18").Value = 2025100723
$ws.Range("F18").Value = 4397
$ws.Range("G18").Value = "三重仁旺店"
$ws.Range("H18").Value = "新北市三重區"
$ws.Range("Q18").Value = "THILF04397"
$ws.Range("R18").Value = "新北一"
$ws.Range("S18").Value = "吳宗鴻"
$ws.Range("T18").Value = 1
$ws.Range("U18").Value = "已完工"
$ws.Range("V18").Value = "2025-10-03 16:12:37"
$ws.Range("W18").Value = "2025-10-03 15:50:00"
$ws.Range("X18").Value = "2025-10-03 16:12:00"
$ws.Range("Z18").Value = 0.4
$ws.Range("AB18").Value = "到場處理"
$ws.Range("AC18").Value = "PMQ4"
$ws.Range("AD18").Value = "O"
$ws.Range("AK18").Value = "O"

# --- 7. Row 19 (item 17) - 三重三文店 / THILF04586 ---
$ws.Range("A19").Value = 17
$ws.Range("B19").Value = "服務"
$ws.Range("C19").Value = 2025100731
$ws.Range("F19").Value = 4586
$ws.Range("G19").Value = "三重三文店"
$ws.Range("H19").Value = "新北市三重區"
$ws.Range("Q19").Value = "THILF04586"
$ws.Range("R19").Value = "新北一"
$ws.Range("S19").Value = "吳宗鴻"
$ws.Range("T19").Value = 1
$ws.Range("U19").Value = "已完工"
$ws.Range("V19").Value = "2025-10-03 16:40:38"
$ws.Range("W19").Value = "2025-10-03 16:20:00"
$ws.Range("X19").Value = "2025-10-03 16:39:00"
$ws.Range("Z19").Value = 0.3
$ws.Range("AB19").Value = "到場處理"
$ws.Range("AC19").Value = "PMQ4"
$ws.Range("AD19").Value = "O"
$ws.Range("AK19").Value = "O"

# --- 8. Print area now covers the extended table, and the saved selection
#        matches where the user last clicked (A19). ---
$ws.PageSetup.PrintArea = '$A$1:$AK$19'
$ws.Range("A19").Select() | Out-Null

"done"
